$d = $word.ActiveDocument

# 1) Drop the leading dialogue tag from the start of the big backstory paragraph
$d.Content.Find.Execute('“Where were you flying?” Shelton lifted', $true, $false, $false, $false, $false, $true, 1, $false, 'Shelton lifted', 2) | Out-Null

# 2) Add one extra trailing space to the end of that paragraph (before its own paragraph mark)
$p = $d.Paragraphs(29)
$r = $p.Range
$endAnchor = $d.Range($r.End - 1, $r.End - 1)
$endAnchor.InsertBefore(' ')

# 3) Open a new paragraph right after it and build out the new content there
$p = $d.Paragraphs(29)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$idx = 30

$para = $d.Paragraphs($idx)
$para.Range.InsertBefore('“Where were you flying?” Shelton asked.')
$para.Range.InsertParagraphAfter()
$idx = $idx + 1

$para = $d.Paragraphs($idx)
$para.Range.InsertBefore('“Denver.” ')
$para.Range.InsertParagraphAfter()
$idx = $idx + 1

$para = $d.Paragraphs($idx)
$para.Range.InsertBefore('On his way back to the sink to wring out the muddy towel he’d been using to scrub the floor, Shelton glanced at the clock set on the kitchen wall – 6:09. “Ah, dammit,” Shelton dropped the towel with a splat into the sink and dashed to the telephone. Nine minutes late. He’d be lucky if April wasn’t already on her way up to the cabin. First the pilot interrupted his shave and now he was giving April a coronary by making Shelton lose track of time. He dialed as quickly as he could and pressed the receiver to his ear.')
$para.Range.InsertParagraphAfter()
$idx = $idx + 1

$para = $d.Paragraphs($idx)
$para.Range.InsertBefore('“Who are you calling?”')
$para.Range.InsertParagraphAfter()
$idx = $idx + 1

$para = $d.Paragraphs($idx)
$para.Range.InsertBefore('“My sister. I check in with her. Let her know I’m okay.”')
$para.Range.InsertParagraphAfter()
$idx = $idx + 1

# three blank paragraphs
for ($k = 0; $k -lt 3; $k++) {
    $para = $d.Paragraphs($idx)
    $para.Range.InsertParagraphAfter()
    $idx = $idx + 1
}

$para = $d.Paragraphs($idx)
$para.Range.InsertBefore('Strengths')
$para.Range.InsertParagraphAfter()
$idx = $idx + 1

$para = $d.Paragraphs($idx)
$para.Range.InsertBefore('Overall, the script is wonderfully written and constructed. I get a strong sense of authorial voice and command. The settings across the Bay Area from sun-baked 70s San Francisco to the streets of Oakland feel well observed and authentic. I see comparisons to existing, successful series - procedural elements mixed with the dissection of serial-killer behavior and motivation of Mindhunter and the time-skipping of True Detective or Bodies. I point this out not to say that it feels derivative, just to say that it suggests there is a market and audience for this type of story. Our link to the present is a true-crime podcaster who gets drawn into the crime she’s reporting on. It’s a solid hook. The crimes touch on sensitive subjects and marginalized people but the story doesn’t come across as condescending or tokenizing. Everything feels real. Throughout this episode, the tension builds nicely leading up to a gripping climax. The structure is generally clear and well-paced. Based upon this episode, I can see intriguing threads of the stories that might be built out throughout the season and the potential for a deep and complex drama. Introducing a political angle through the new mayor of San Francisco will bring added dimensions to the overarching story and themes.     ')
$para.Range.InsertParagraphAfter()
$idx = $idx + 1

# one blank paragraph
$para = $d.Paragraphs($idx)
$para.Range.InsertParagraphAfter()
$idx = $idx + 1

$para = $d.Paragraphs($idx)
$para.Range.InsertBefore('Weaknesses')
$para.Range.InsertParagraphAfter()
$idx = $idx + 1

$para = $d.Paragraphs($idx)
$para.Range.InsertBefore('I think what could make this story even stronger for me is if you examined the procedural aspects of the script and looked for ways to make them stand out more. There are dozens of series that focus on police or investigators doing their jobs. You should be sure that you’re offering something unique to avoid getting lost in the shuffle. We end up meeting a lot of police officers (even one of the victim’s children grows up to be a police officer) and they run together a little. At least in what I’ve read, two of the three main characters we’re introduced to – Seth and Dennis – feel too similar. Maybe this changes in future episodes but more distinctive voices for these two characters would make each of the different time periods equally compelling. I was somewhat thrown off by the introduction to Seth. It seems as if both he and his supervisor are disappointed that he didn’t shoot a fleeing suspect in the back. Isn’t that a good thing? I don’t think it even takes suffering a trauma like the one Seth did to feel like shooting a suspect over a drug offense is unnecessary. A little bit of a disconnect there. These are minor adjustments that could make a strong script even better.')

Write-Output "done"
